$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.478.10"
$ws.Range("E2").Value = "  +2.02%  "
$ws.Range("D3").Value = "2.327.60"
$ws.Range("E3").Value = "  +1.49%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'545.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").Value = "'132.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'0.585"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.20%  "
$ws.Range("D9").Value = "2.318.27"
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("D10").Value = "'0.100"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.79%  "
$ws.Range("D11").Value = "'5.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").Value = "'0.334"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("D14").Value = "'23.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("D15").Value = "2.736.57"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("D16").Value = "59.288.11"
$ws.Range("E16").Value = "  +1.78%  "
$ws.Range("D17").Value = "'0.0000133"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.87%  "
$ws.Range("D18").Value = "2.310.13"
$ws.Range("E18").Value = "  +3.55%  "
$ws.Range("D19").Value = "'10.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "'4.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("D21").Value = "'314.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("D22").Value = "'6.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.30%  "
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "'62.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "'0.173"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.10%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").Value = "'7.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Value = "'1.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.10%  "
$ws.Range("D31").Value = "'1.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.17%  "
$ws.Range("D32").Value = "0.0₃0741"
$ws.Range("E32").Value = "  +3.26%  "
$ws.Range("D33").Value = "'5.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.85%  "
$ws.Range("D34").Value = "'1.42"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +14.96%  "
$ws.Range("D35").Value = "'0.385"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.64%  "
$ws.Range("D37").Value = "'17.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.53%  "
$ws.Range("D38").Value = "'0.996"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("D39").Value = "'4.08"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.12%  "
$ws.Range("D40").Value = "'318.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.98%  "
$ws.Range("D41").Value = "'38.04"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").Value = "'1.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.52%  "
$ws.Range("D43").Value = "'143.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.58%  "
$ws.Range("D44").Value = "'3.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.70%  "
$ws.Range("D45").Value = "'0.0953"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("D46").Value = "'0.0496"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.75%  "
$ws.Range("D47").Value = "'0.559"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.24%  "
$ws.Range("D48").Value = "'18.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.04%  "
$ws.Range("D49").Value = "'0.0210"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("E51").Value = "  -0.07%  "

# Rows 29 and 30 swap: PancakeSwap moves to row 29, Monero moves to row 30
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'1.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "'171.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.72%  "
